$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$b2 = $ws.Range("B2").Value2
$ws.Range("B2").Value2 = $b2 + "`n order By ss.study_subject_id ASC LIMIT 100"

$b3 = $ws.Range("B3").Value2
$ws.Range("B3").Value2 = $b3 + "`n order By samp.sample_id ASC LIMIT 100"

$b4 = $ws.Range("B4").Value2
$ws.Range("B4").Value2 = $b4 + "`n order By f.file_name ASC LIMIT 100"

[void]$ws.Range("B2").Select()
